# Daily attendance processing - 2025-11-30 10:26:29
#
# The automated attendance-processing run re-evaluated session statuses and
# recorded-by lists. Net effect on the sheet:
#   - Several "Recorded By" cells (G2, G3, G4, G5, G7, G12, G28) were
#     rewritten with the same set of reporters in a different (re-sorted)
#     order.
#   - The PARASITOLOGY session on 30/11/2025 (row 16) moved from "Pending"
#     to "Not Recorded" now that its date has passed with 0 recorded
#     attendance - its status text and row shading (pending-yellow ->
#     not-recorded-pink, same shading already used by row 29) changed
#     together.
#   - The summary counters that track this (Missing/Pending Sessions on
#     the dashboard, and the matching Year2/C1 row in the pivot table)
#     were updated to reflect the one-session shift from Pending to
#     Missing/Not Recorded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Recorded By cells: same reporters, refreshed ordering -----------------
$ws.Range("G2").Value = "Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, System, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G3").Value = "majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, System"
$ws.Range("G4").Value = "majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G5").Value = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G7").Value = "lamiaa.ossama@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg"
$ws.Range("G12").Value = "Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
$ws.Range("G28").Value = "maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"

# --- Row 16 (PARASITOLOGY, 30/11/2025): Pending -> Not Recorded ------------
# Re-use row 29's existing "Not Recorded" formatting (same pink shading)
# instead of inventing a new fill, then flip the status text.
$ws.Range("A29:I29").Copy()
$ws.Range("A16:I16").PasteSpecial(-4122)
$ws.Range("I16").Value = "Not Recorded"

# --- Dashboard counters: one session shifted Pending -> Missing ------------
$ws.Range("L7").Value = 2
$ws.Range("L8").Value = 11

# --- Pivot-style summary row (Year2/C1): same shift -------------------------
$ws.Range("P15").Value = 2
$ws.Range("Q15").Value = 11
